# Add a new "Bat" (battery/appliance charge level) column to the
# explicaciones sheet: header in I1, Spanish text in I2, English text in I3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "Bat"
$ws.Range("I3").Value = "Appliance charge level"
$ws.Range("I2").Value = "Nivel de carga de los aparatos"

$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("J9").Select()
